# ---------------------------------------------------------------------------
# Edit: swap the presentation's theme color scheme from the "Integral" theme
# to the "Office Theme" colors, and re-style the plenary table on slide 16
# with the new (built-in) table style.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" table)
#    {A74A870B-8196-46C4-8386-E778ADFAA569} -> {95147D40-2931-4D54-B2A8-1E2E544FE72C}
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle('{95147D40-2931-4D54-B2A8-1E2E544FE72C}')
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colours: replace the "Integral" palette with the stock
#    "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$newColors = @(
    @(0,0,0),          # 1  dk1
    @(255,255,255),    # 2  lt1
    @(68,84,106),       # 3  dk2
    @(231,230,230),     # 4  lt2
    @(91,155,213),      # 5  accent1
    @(237,125,49),      # 6  accent2
    @(165,165,165),     # 7  accent3
    @(255,192,0),       # 8  accent4
    @(68,114,196),      # 9  accent5
    @(112,173,71),      # 10 accent6
    @(5,99,193),        # 11 hlink
    @(149,79,114)       # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $newColors[$i - 1]
    $packed = $rgb[0] + ($rgb[1] * 256) + ($rgb[2] * 65536)
    $colorScheme.Item($i).RGB = $packed
}
